# Auto-generated edit script applying the Typhon_Profits leve-price refresh
# (H=currentAveragePrice, I=currentAveragePriceNQ, J=currentAveragePriceHQ,
#  K=LevePriceNQ, L=LevePriceHQ, M=LeveProfitNQ, N=LeveProfitHQ)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 3161
$ws.Range("I62").Value = 3301.6667
$ws.Range("J62").Value = 2950
$ws.Range("K62").Value = 3301.6667
$ws.Range("L62").Value = 2950
$ws.Range("M62").Value = -2677.6667
$ws.Range("N62").Value = -4198
$ws.Range("H65").Value = 3161
$ws.Range("I65").Value = 3301.6667
$ws.Range("J65").Value = 2950
$ws.Range("K65").Value = 16508.3335
$ws.Range("L65").Value = 14750
$ws.Range("M65").Value = -13388.3335
$ws.Range("N65").Value = -20990
$ws.Range("H129").Value = 848.175
$ws.Range("J129").Value = 850.7105
$ws.Range("L129").Value = 2552.1315
$ws.Range("N129").Value = -12552.1315
$ws.Range("H132").Value = 3979.0454
$ws.Range("I132").Value = 4031.0557
$ws.Range("K132").Value = 12093.1671
$ws.Range("M132").Value = -9563.167099999999
$ws.Range("H137").Value = 61606.59
$ws.Range("I137").Value = 3983.1667
$ws.Range("K137").Value = 11949.5001
$ws.Range("M137").Value = -9399.500100000001
$ws.Range("H140").Value = 55555
$ws.Range("J140").Value = 55555
$ws.Range("L140").Value = 55555
$ws.Range("N140").Value = -65915
$ws.Range("H141").Value = 2473.25
$ws.Range("I141").Value = 1742.2222
$ws.Range("K141").Value = 5226.6666
$ws.Range("M141").Value = -46.66659999999956

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 16531.582
$ws.Range("I32").Value = 17329.08
$ws.Range("K32").Value = 17329.08
$ws.Range("M32").Value = -17042.08
$ws.Range("H45").Value = 3173.1667
$ws.Range("I45").Value = 2724.9
$ws.Range("J45").Value = 3493.3572
$ws.Range("K45").Value = 2724.9
$ws.Range("L45").Value = 3493.3572
$ws.Range("M45").Value = -2347.9
$ws.Range("N45").Value = -4247.3572
$ws.Range("H61").Value = 3409.389
$ws.Range("I61").Value = 2459.1538
$ws.Range("J61").Value = 5880
$ws.Range("K61").Value = 2459.1538
$ws.Range("L61").Value = 5880
$ws.Range("M61").Value = -2247.1538
$ws.Range("N61").Value = -6304
$ws.Range("H74").Value = 2136.4443
$ws.Range("I74").Value = 2194.9565
$ws.Range("J74").Value = 1800
$ws.Range("K74").Value = 2194.9565
$ws.Range("L74").Value = 1800
$ws.Range("M74").Value = -1320.9565
$ws.Range("N74").Value = -3548
$ws.Range("H77").Value = 2136.4443
$ws.Range("I77").Value = 2194.9565
$ws.Range("J77").Value = 1800
$ws.Range("K77").Value = 10974.7825
$ws.Range("L77").Value = 9000
$ws.Range("M77").Value = -6606.782499999999
$ws.Range("N77").Value = -17736
$ws.Range("H97").Value = 1862
$ws.Range("I97").Value = 1744.1538
$ws.Range("J97").Value = 2245
$ws.Range("K97").Value = 1744.1538
$ws.Range("L97").Value = 2245
$ws.Range("M97").Value = -1248.1538
$ws.Range("N97").Value = -3237
$ws.Range("H122").Value = 1291.45
$ws.Range("I122").Value = 1417.9333
$ws.Range("J122").Value = 912
$ws.Range("K122").Value = 4253.7999
$ws.Range("L122").Value = 2736
$ws.Range("M122").Value = -1803.7999
$ws.Range("N122").Value = -7636
$ws.Range("H132").Value = 24437.61
$ws.Range("I132").Value = 2565.6155
$ws.Range("J132").Value = 52871.2
$ws.Range("K132").Value = 7696.8465
$ws.Range("L132").Value = 158613.6
$ws.Range("M132").Value = -5166.8465
$ws.Range("N132").Value = -163673.6
$ws.Range("H136").Value = 3409.389
$ws.Range("I136").Value = 2459.1538
$ws.Range("J136").Value = 5880
$ws.Range("K136").Value = 7377.4614
$ws.Range("L136").Value = 17640
$ws.Range("M136").Value = -4827.4614
$ws.Range("N136").Value = -22740

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 921
$ws.Range("I20").Value = 947
$ws.Range("J20").Value = 863.8
$ws.Range("K20").Value = 947
$ws.Range("L20").Value = 863.8
$ws.Range("M20").Value = -700
$ws.Range("N20").Value = -1357.8
$ws.Range("H22").Value = 204
$ws.Range("I22").Value = 204
$ws.Range("K22").Value = 204
$ws.Range("M22").Value = -31
$ws.Range("H134").Value = 41992
$ws.Range("I134").Value = 54039.7
$ws.Range("K134").Value = 162119.1
$ws.Range("M134").Value = -159584.1

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 820
$ws.Range("I22").Value = 250
$ws.Range("J22").Value = 1390
$ws.Range("K22").Value = 250
$ws.Range("L22").Value = 1390
$ws.Range("M22").Value = 100
$ws.Range("N22").Value = -2090
$ws.Range("H31").Value = 1714.6364
$ws.Range("I31").Value = 1148.5
$ws.Range("J31").Value = 4262.25
$ws.Range("K31").Value = 1148.5
$ws.Range("L31").Value = 4262.25
$ws.Range("M31").Value = -853.5
$ws.Range("N31").Value = -4852.25
$ws.Range("H34").Value = 1714.6364
$ws.Range("I34").Value = 1148.5
$ws.Range("J34").Value = 4262.25
$ws.Range("K34").Value = 1148.5
$ws.Range("L34").Value = 4262.25
$ws.Range("M34").Value = -946.5
$ws.Range("N34").Value = -4666.25
$ws.Range("H70").Value = 16045
$ws.Range("J70").Value = 16045
$ws.Range("L70").Value = 16045
$ws.Range("N70").Value = -16675
$ws.Range("H73").Value = 16045
$ws.Range("J73").Value = 16045
$ws.Range("L73").Value = 16045
$ws.Range("N73").Value = -18229
$ws.Range("H122").Value = 977
$ws.Range("I122").Value = 976.5
$ws.Range("J122").Value = 980
$ws.Range("K122").Value = 2929.5
$ws.Range("L122").Value = 2940
$ws.Range("M122").Value = -479.5
$ws.Range("N122").Value = -7840

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H109").Value = 4908.143
$ws.Range("I109").Value = 1203.3334
$ws.Range("J109").Value = 5918.5454
$ws.Range("K109").Value = 3610.0002
$ws.Range("L109").Value = 17755.6362
$ws.Range("M109").Value = -2570.0002
$ws.Range("N109").Value = -19835.6362

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("N93").ClearContents()
$ws.Range("H97").Value = 2564.4783
$ws.Range("J97").Value = 4863.2856
$ws.Range("L97").Value = 4863.2856
$ws.Range("N97").Value = -5855.2856
$ws.Range("H113").Value = 2876.6667
$ws.Range("I113").Value = 1997.5
$ws.Range("J113").Value = 3580
$ws.Range("K113").Value = 1997.5
$ws.Range("L113").Value = 3580
$ws.Range("M113").Value = 172.5
$ws.Range("N113").Value = -7920
$ws.Range("H122").Value = 1843.2142
$ws.Range("I122").Value = 1888.7778
$ws.Range("J122").Value = 1761.2
$ws.Range("K122").Value = 5666.3334
$ws.Range("L122").Value = 5283.6
$ws.Range("M122").Value = -3216.3334
$ws.Range("N122").Value = -10183.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4545.364
$ws.Range("I7").Value = 4662.375
$ws.Range("K7").Value = 4662.375
$ws.Range("M7").Value = -4550.375
$ws.Range("H40").Value = 5492.7856
$ws.Range("I40").Value = 5056.2856
$ws.Range("J40").Value = 5929.2856
$ws.Range("K40").Value = 5056.2856
$ws.Range("L40").Value = 5929.2856
$ws.Range("M40").Value = -4920.2856
$ws.Range("N40").Value = -6201.2856
$ws.Range("H122").Value = 1035373.2
$ws.Range("I122").Value = 2181227.2
$ws.Range("J122").Value = 4104.4
$ws.Range("K122").Value = 6543681.600000001
$ws.Range("L122").Value = 12313.2
$ws.Range("M122").Value = -6541231.600000001
$ws.Range("N122").Value = -17213.2
$ws.Range("H126").Value = 4545.364
$ws.Range("I126").Value = 4662.375
$ws.Range("K126").Value = 13987.125
$ws.Range("M126").Value = -11517.125
$ws.Range("H132").Value = 2368.2666
$ws.Range("I132").Value = 1692.8182
$ws.Range("K132").Value = 5078.4546
$ws.Range("M132").Value = -2548.4546

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 5683535.5
$ws.Range("I107").Value = 579.6667
$ws.Range("K107").Value = 1739.0001
$ws.Range("M107").Value = 180.9999
$ws.Range("H119").Value = 29999.715
$ws.Range("J119").Value = 29999.715
$ws.Range("L119").Value = 29999.715
$ws.Range("N119").Value = -39675.715
$ws.Range("H122").Value = 1754.409
$ws.Range("J122").Value = 1823.8334
$ws.Range("L122").Value = 5471.5002
$ws.Range("N122").Value = -10371.5002
$ws.Range("H136").Value = 22728742
$ws.Range("I136").Value = 38463030
$ws.Range("J136").Value = 1433.7778
$ws.Range("K136").Value = 115389090
$ws.Range("L136").Value = 4301.3334
$ws.Range("M136").Value = -115386540

